# Dream Checklist.docx edit script
# Summary of changes (per commit "Can parse tags on creation."):
#  1. Remove the "Text posts" To-Do bullet entirely.
#  2. Move the hidden "_GoBack" bookmark from the end of the document
#     (after "...drag to upload type thing") onto the "Front-End" bullet.
#  3. Tighten "Formatting the captions and tags" -> "Formatting the tags".
#  4. Fix "Autoplay" -> "Auto-play" in the autoplay videos/gifs bullet.
#  5. Add a new Bugs bullet "Tilted Samsung pictures" right after
#     "Does hotness rating actually work?".

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, [string]$needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# --- 1. Delete the "Text posts" paragraph -------------------------------
$textPosts = Find-ParagraphByText $d "Text posts"
if ($textPosts -ne $null) {
    $textPosts.Range.Delete()
}

# --- 2. Move the _GoBack bookmark to the "Front-End" paragraph ----------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$frontEnd = Find-ParagraphByText $d "Front-End"
$frontEndStart = $d.Range($frontEnd.Range.Start, $frontEnd.Range.Start)
$d.Bookmarks.Add("_GoBack", $frontEndStart)

# --- 3. "Formatting the captions and tags" -> "Formatting the tags" -----
$d.Content.Find.Execute("captions and ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- 4. "Autoplay videos..." -> "Auto-play videos..." -------------------
$d.Content.Find.Execute("Autoplay videos", $true, $false, $false, $false, $false, $true, 1, $false, "Auto-play videos", 2)

# --- 5. Insert new Bugs bullet "Tilted Samsung pictures" ----------------
$hotness = Find-ParagraphByText $d "Does hotness rating actually work?"
$hotness.Range.InsertParagraphAfter()
$tilted = $hotness.Next()
$tilted.Range.Text = "Tilted Samsung pictures"
$tilted.Range.ListFormat.ListLevelNumber = 3
